$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '22.965.04'
$ws.Range("E2").Value = '  -3.96%  '
$ws.Range("D3").Value = '1.599.58'
$ws.Range("E3").Value = '  -2.96%  '
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").Value = '1.002'
$ws.Range("E5").Value = '  +0.21%  '
$ws.Range("D6").Value = "'" + '300.30'
$ws.Range("E6").Value = '  -3.33%  '
$ws.Range("D7").Value = '0.3769'
$ws.Range("E7").Value = '  -3.21%  '
$ws.Range("D8").Value = '0.3617'
$ws.Range("E8").Value = '  -5.72%  '
$ws.Range("D9").Value = '49.56'
$ws.Range("E9").Value = '  -2.98%  '
$ws.Range("D10").Value = '1.249'
$ws.Range("E10").Value = '  -6.90%  '
$ws.Range("D11").Value = '1.003'
$ws.Range("E11").Value = '  +0.29%  '
$ws.Range("D12").Value = '0.08083'
$ws.Range("E12").Value = '  -4.22%  '
$ws.Range("D13").Value = '22.42'
$ws.Range("E13").Value = '  -6.19%  '
$ws.Range("D14").Value = '6.555'
$ws.Range("E14").Value = '  -6.59%  '
$ws.Range("D15").Value = '7.308'
$ws.Range("E15").Value = '  -7.11%  '
$ws.Range("D16").Value = "'" + '0.00001230'
$ws.Range("E16").Value = '  -6.70%  '
$ws.Range("D17").Value = '1.595.46'
$ws.Range("E17").Value = '  -3.23%  '
$ws.Range("D18").Value = '91.85'
$ws.Range("E18").Value = '  -2.29%  '
$ws.Range("D19").Value = '0.06845'
$ws.Range("E19").Value = '  -1.86%  '
$ws.Range("D20").Value = '18.09'
$ws.Range("E20").Value = '  -7.69%  '
$ws.Range("D21").Value = '6.527'
$ws.Range("E21").Value = '  -5.73%  '
$ws.Range("D22").Value = '0.5576'
$ws.Range("E22").Value = '  -5.55%  '
$ws.Range("D23").Value = '1.003'
$ws.Range("E23").Value = '  +0.26%  '
$ws.Range("D24").Value = '13.04'
$ws.Range("E24").Value = '  -4.55%  '
$ws.Range("D25").Value = '22.965.72'
$ws.Range("E25").Value = '  -4.00%  '
$ws.Range("D26").Value = '2.359'
$ws.Range("E26").Value = '  -3.28%  '
$ws.Range("D27").Value = '2.795'
$ws.Range("E27").Value = '  -4.78%  '
$ws.Range("D28").Value = '21.02'
$ws.Range("E28").Value = '  -4.23%  '
$ws.Range("D29").Value = "'" + '150.10'
$ws.Range("E29").Value = '  -2.48%  '
$ws.Range("D30").Value = '5.255'
$ws.Range("E30").Value = '  -3.57%  '
$ws.Range("D31").Value = '133.05'
$ws.Range("E31").Value = '  -2.94%  '
$ws.Range("D32").Value = '2.275'
$ws.Range("E32").Value = '  -9.12%  '
$ws.Range("D33").Value = '6.708'
$ws.Range("E33").Value = '  -13.12%  '
$ws.Range("D34").Value = '1.775.98'
$ws.Range("E34").Value = '  -3.06%  '
$ws.Range("D35").Value = '0.9527'
$ws.Range("E35").Value = '  -3.45%  '
$ws.Range("D36").Value = '0.07534'
$ws.Range("E36").Value = '  -7.20%  '
$ws.Range("D37").Value = '10.27'
$ws.Range("E37").Value = '  -1.86%  '
$ws.Range("D38").Value = '6.176'
$ws.Range("E38").Value = '  -7.91%  '
$ws.Range("B39").Value = 'Algorand'
$ws.Range("C39").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D39").Value = '0.2512'
$ws.Range("E39").Value = '  -6.58%  '
$ws.Range("B40").Value = 'Stellar'
$ws.Range("C40").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D40").Value = '0.08829'
$ws.Range("E40").Value = '  -3.04%  '
$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").Value = '0.02672'
$ws.Range("E41").Value = '  -8.81%  '
$ws.Range("D42").Value = '1.367'
$ws.Range("E42").Value = '  -4.04%  '
$ws.Range("D43").Value = '0.6965'
$ws.Range("E43").Value = '  -7.72%  '
$ws.Range("D44").Value = '12.25'
$ws.Range("E44").Value = '  -8.55%  '
$ws.Range("E45").Value = '  -8.15%  '
$ws.Range("D46").Value = '0.6551'
$ws.Range("E46").Value = '  -5.51%  '
$ws.Range("D47").Value = '1.001'
$ws.Range("E47").Value = '  +0.17%  '
$ws.Range("D48").Value = '2.286'
$ws.Range("E48").Value = '  -6.22%  '
$ws.Range("D49").Value = '3.978'
$ws.Range("E49").Value = '  -2.89%  '
$ws.Range("D50").Value = '131.84'
$ws.Range("E50").Value = '  -1.93%  '
$ws.Range("D51").Value = '0.07895'
